$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column D, row 2: formula C2^-1
$ws.Range("D2").Formula = "=C2^-1"
$ws.Range("D2").NumberFormat = $ws.Range("C2").NumberFormat

# Update selection to match the authored state
[void]$ws.Range("I13").Select()
